$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("freq")

# Updated models for experiment 1: rename the "exp" column values.
# D2:D81 were "lago" -> "lagoetal"
$ws.Range("D2:D81").Value = "lagoetal"
# D82:D161 were "up" -> "exp1"
$ws.Range("D82:D161").Value = "exp1"

# Changed the text location: move the view/selection from B40 to K90,
# scrolling the visible top-left cell from A36 to A73.
$ws.Range("K90").Select()
$win = $excel.ActiveWindow()
$win.ScrollRow = 73
$win.ScrollColumn = 1
